# Weekly update: insert a new Lechuga price record for week of 2021-11-09
# (serial 44509) into the "Vega Modelo de Temuco" sheet, just before the
# existing 44421 records, pushing everything from the old row 757 onward
# down by one row (785 -> 786).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 757; Excel shifts rows 757:785 down to 758:786.
$ws.Rows.Item(757).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A757").Value = 10
$ws.Range("B757").Value = "Vega Modelo de Temuco"
$ws.Range("C757").Value = "La Araucanía"
$ws.Range("D757").Value = 44509
$ws.Range("E757").Value = 9
$ws.Range("F757").Value = 100112033
$ws.Range("G757").Value = "Lechuga"
$ws.Range("H757").Value = "Escarola"
$ws.Range("I757").Value = "Primera"
$ws.Range("J757").Value = 800
$ws.Range("K757").Value = 7000
$ws.Range("L757").Value = 7000
$ws.Range("M757").Value = 7000
$ws.Range("N757").Value = "`$/caja 15 unidades"
$ws.Range("O757").Value = "Región del Maule"
$ws.Range("P757").Value = 467
$ws.Range("Q757").Value = 15
$ws.Range("R757").Value = "Hortaliza"
